$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Create rows 7 and 8 first (brand-new rows/strings) before touching existing B2:B6
# so the new "ANN(2HL-19,18 Nodes)" / "ANN(2HL-19,20 Nodes)" strings are appended
# cleanly while the old model-name strings are still present/referenced.

# --- Row 7 (new) ---
$ws.Range("A7").Value = 55
$ws.Range("B7").Value = "ANN(2HL-19,18 Nodes)"
$ws.Range("C7").Value = 14.43398277484873
$ws.Range("D7").Value = 0.02001350559100407
$ws.Range("E7").Value = 0.2052919306225131
$ws.Range("F7").Value = 24.99954061411335
$ws.Range("G7").Value = 0.9558272651872618
$ws.Range("H7").Value = 0.8305331942045113
$ws.Range("I7").Value = 0.8329472191054043
$ws.Range("J7").Value = 0.002906596530684237
$ws.Range("K7").Value = 0.8566798373674943
$ws.Range("L7").Value = 0.8590529495708403
$ws.Range("M7").Value = 0.002770127298242939

# --- Row 8 (new) ---
$ws.Range("A8").Value = 57
$ws.Range("B8").Value = "ANN(2HL-19,20 Nodes)"
$ws.Range("C8").Value = 14.04530482976033
$ws.Range("D8").Value = 0.01732375264552129
$ws.Range("E8").Value = 0.1935434560399775
$ws.Range("F8").Value = 24.32640549076017
$ws.Range("G8").Value = 0.9617869969330447
$ws.Range("H8").Value = 0.8305331942045113
$ws.Range("I8").Value = 0.8297178064735136
$ws.Range("J8").Value = 0.0009817641687142262
$ws.Range("K8").Value = 0.8566798373674943
$ws.Range("L8").Value = 0.8563907633913891
$ws.Range("M8").Value = 0.0003374352511826652

# --- Row 2 (existing, update in place) ---
$ws.Range("A2").Value = 3
$ws.Range("B2").Value = "ANN(2HL-14,20 Nodes)"
$ws.Range("C2").Value = 14.97766183632132
$ws.Range("D2").Value = 0.0196565875068293
$ws.Range("E2").Value = 0.1994921902261521
$ws.Range("F2").Value = 25.94129678204609
$ws.Range("G2").Value = 0.9568108120973831
$ws.Range("H2").Value = 0.8305331942045113
$ws.Range("I2").Value = 0.8324456631738312
$ws.Range("J2").Value = 0.002302700220370625
$ws.Range("K2").Value = 0.8566798373674943
$ws.Range("L2").Value = 0.8621402850447035
$ws.Range("M2").Value = 0.006373965440798411

# --- Row 3 (existing, update in place) ---
$ws.Range("A3").Value = 13
$ws.Range("B3").Value = "ANN(2HL-15,14 Nodes)"
$ws.Range("C3").Value = 14.90953999469957
$ws.Range("D3").Value = 0.02519401325951733
$ws.Range("E3").Value = 0.1983441992264328
$ws.Range("F3").Value = 25.82330678671821
$ws.Range("G3").Value = 0.9500732051899409
$ws.Range("H3").Value = 0.8305331942045113
$ws.Range("I3").Value = 0.8307020771951377
$ws.Range("J3").Value = 0.0002033428546923967
$ws.Range("K3").Value = 0.8566798373674943
$ws.Range("L3").Value = 0.8619133396454793
$ws.Range("M3").Value = 0.006109052705228993

# --- Row 4 (existing, update in place) ---
$ws.Range("A4").Value = 24
$ws.Range("B4").Value = "ANN(2HL-17,16 Nodes)"
$ws.Range("C4").Value = 14.11485799363989
$ws.Range("D4").Value = 0.01955891295117575
$ws.Range("E4").Value = 0.204703539792901
$ws.Range("F4").Value = 24.44678634200373
$ws.Range("G4").Value = 0.956884842150962
$ws.Range("H4").Value = 0.8305331942045113
$ws.Range("I4").Value = 0.8314579319455953
$ws.Range("J4").Value = 0.001113426588529846
$ws.Range("K4").Value = 0.8566798373674943
$ws.Range("L4").Value = 0.8613865444695308
$ws.Range("M4").Value = 0.005494126156277689

# --- Row 5 (existing, update in place) ---
$ws.Range("A5").Value = 26
$ws.Range("B5").Value = "ANN(2HL-17,18 Nodes)"
$ws.Range("C5").Value = 14.00944728863614
$ws.Range("D5").Value = 0.0230771304572704
$ws.Range("E5").Value = 0.2119809145174359
$ws.Range("F5").Value = 24.26413756013499
$ws.Range("G5").Value = 0.9510604696987351
$ws.Range("H5").Value = 0.8305331942045113
$ws.Range("I5").Value = 0.8333274768543103
$ws.Range("J5").Value = 0.003364444274229582
$ws.Range("K5").Value = 0.8566798373674943
$ws.Range("L5").Value = 0.8562385170549937
$ws.Range("M5").Value = 0.0005151519777292761

# --- Row 6 (existing, update in place) ---
$ws.Range("A6").Value = 51
$ws.Range("B6").Value = "ANN(2HL-19,16 Nodes)"
$ws.Range("C6").Value = 14.07707741717936
$ws.Range("D6").Value = 0.02070188394352684
$ws.Range("E6").Value = 0.1985064035219516
$ws.Range("F6").Value = 24.3813964421143
$ws.Range("G6").Value = 0.9570052578144018
$ws.Range("H6").Value = 0.8305331942045113
$ws.Range("I6").Value = 0.8283210172955691
$ws.Range("J6").Value = 0.002663562304768684
$ws.Range("K6").Value = 0.8566798373674943
$ws.Range("L6").Value = 0.8614416923324251
$ws.Range("M6").Value = 0.005558500103800201
